$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear existing used range first
$ws.Cells.Clear()

$data = @(
    @('Headline Title', 'URL'),
    @('France starts collecting tax on tech giants', 'http://techcrunch.com/2020/11/25/france-starts-collecting-tax-on-tech-giants/'),
    @('France to Tech Giants: Pay Up, Buttercup', 'https://gizmodo.com/france-to-tech-giants-pay-up-buttercup-1845757236'),
    @('Why France may ban discrimination against accents', 'https://www.bbc.co.uk/news/world-europe-55069048'),
    @('Coronavirus: France replaces lockdown with evening curfew', 'https://www.bbc.co.uk/news/world-europe-55324422'),
    @('France police security bill: Protests turn violent again', 'https://www.bbc.co.uk/news/world-europe-55201993'),
    @('Where Are Your Papers? A New Yorker Goes Home to Locked-Down France', 'https://www.nytimes.com/2020/12/10/us/france-lockdown-compare-new-york.html'),
    @('Patrick Vieira: Nice sack manager after two and a half years in charge', 'https://www.bbc.co.uk/sport/football/55187481'),
    @('Covid: Students and retirees form long-distance friendships', 'https://www.bbc.co.uk/news/world-europe-55235378'),
    @('Charlie Hebdo: Fourteen guilty in 2015 Paris terror attacks trial', 'https://www.bbc.co.uk/news/world-europe-55336094'),
    @('Trial of ex-president Sarkozy a landmark for France', 'https://www.bbc.co.uk/news/world-europe-55015479'),
    @('France Islam: Muslims face state pressure to embrace values', 'https://www.bbc.co.uk/news/world-europe-55132098'),
    @('Radio France Internationale publishes obituaries of people still alive', 'https://www.bbc.co.uk/news/world-europe-54965098'),
    @('France Islam: French cabinet backs law targeting extremism', 'https://www.bbc.co.uk/news/world-europe-55253832'),
    @('Lockdowns Aren’t Funny. But You’ve Still Got to Laugh.', 'https://www.nytimes.com/2020/11/19/theater/france-lockdown-comedy.html'),
    @('Obituary: Valéry Giscard DEstaing', 'https://www.bbc.co.uk/news/world-europe-13062449'),
    @('French police seize 10 tigers after mistreatment complaint', 'https://www.bbc.co.uk/news/world-europe-55337455'),
    @('French rugby legend Dominici dies in park tragedy', 'https://www.bbc.co.uk/news/55056643'),
    @('Ireland and Scotland in same 2023 World Cup pool as South Africa', 'https://www.bbc.co.uk/sport/rugby-union/55301480'),
    @('Man Citys Torres hits hat-trick as Spain put six past Germany to reach Nations League finals', 'https://www.bbc.co.uk/sport/football/54979319'),
    @('Wasps hooker Barbeary called up to England squad', 'https://www.bbc.co.uk/sport/rugby-union/55139252')
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
}
